# Scheduled market-price data refresh
# Updates currentAveragePrice / NQ / HQ price & profit columns (H:N)
# for the leve rows whose cached market data changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 173679.33
$ws.Range("I62").Value = 189250.27
$ws.Range("K62").Value = 189250.27
$ws.Range("M62").Value = -188626.27
# Row 65
$ws.Range("H65").Value = 173679.33
$ws.Range("I65").Value = 189250.27
$ws.Range("K65").Value = 946251.35
$ws.Range("M65").Value = -943131.35
# Row 74
$ws.Range("H74").Value = 5490
$ws.Range("I74").Value = 5490
$ws.Range("K74").Value = 5490
$ws.Range("M74").Value = -4554
# Row 77
$ws.Range("H77").Value = 5490
$ws.Range("I77").Value = 5490
$ws.Range("K77").Value = 27450
$ws.Range("M77").Value = -22770
# Row 86
$ws.Range("H86").Value = 1761.55
$ws.Range("I86").Value = 1890.2354
$ws.Range("K86").Value = 1890.2354
$ws.Range("M86").Value = -767.2354
# Row 89
$ws.Range("H89").Value = 1761.55
$ws.Range("I89").Value = 1890.2354
$ws.Range("K89").Value = 9451.177
$ws.Range("M89").Value = -3835.177
# Row 103
$ws.Range("H103").Value = 583.4737
$ws.Range("J103").Value = 501.23077
$ws.Range("L103").Value = 1503.69231
$ws.Range("N103").Value = -2675.69231
# Row 137
$ws.Range("H137").Value = 7906.4136
$ws.Range("I137").Value = 6251.6772
$ws.Range("K137").Value = 18755.0316
$ws.Range("M137").Value = -16205.0316
# Row 138
$ws.Range("H138").Value = 2150.8
$ws.Range("I138").Value = 1608.5
$ws.Range("J138").Value = 3114.889
$ws.Range("K138").Value = 4825.5
$ws.Range("L138").Value = 9344.667000000001
$ws.Range("M138").Value = 314.5
$ws.Range("N138").Value = -19624.667

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 7872.6665
$ws.Range("I61").Value = 4919.5713
$ws.Range("J61").Value = 12007
$ws.Range("K61").Value = 4919.5713
$ws.Range("L61").Value = 12007
$ws.Range("M61").Value = -4707.5713
$ws.Range("N61").Value = -12431
# Row 74
$ws.Range("H74").Value = 2076.0645
$ws.Range("I74").Value = 912.1539
$ws.Range("J74").Value = 8128.4
$ws.Range("K74").Value = 912.1539
$ws.Range("L74").Value = 8128.4
$ws.Range("M74").Value = -38.15390000000002
$ws.Range("N74").Value = -9876.4
# Row 77
$ws.Range("H77").Value = 2076.0645
$ws.Range("I77").Value = 912.1539
$ws.Range("J77").Value = 8128.4
$ws.Range("K77").Value = 4560.7695
$ws.Range("L77").Value = 40642
$ws.Range("M77").Value = -192.7695000000003
$ws.Range("N77").Value = -49378
# Row 132
$ws.Range("H132").Value = 8871.053
$ws.Range("I132").Value = 7073.2046
$ws.Range("K132").Value = 21219.6138
$ws.Range("M132").Value = -18689.6138
# Row 136
$ws.Range("H136").Value = 7872.6665
$ws.Range("I136").Value = 4919.5713
$ws.Range("J136").Value = 12007
$ws.Range("K136").Value = 14758.7139
$ws.Range("L136").Value = 36021
$ws.Range("M136").Value = -12208.7139
$ws.Range("N136").Value = -41121

$ws = $wb.Worksheets.Item("BSM")
# Row 30
$ws.Range("H30").Value = 6389.8
$ws.Range("J30").Value = 6737.25
$ws.Range("L30").Value = 6737.25
$ws.Range("N30").Value = -6987.25
# Row 64
$ws.Range("H64").Value = 2752
$ws.Range("I64").Value = 2000.5
$ws.Range("J64").Value = 3503.5
$ws.Range("K64").Value = 2000.5
$ws.Range("L64").Value = 3503.5
$ws.Range("M64").Value = -1775.5
$ws.Range("N64").Value = -3953.5
# Row 67
$ws.Range("H67").Value = 2752
$ws.Range("I67").Value = 2000.5
$ws.Range("J67").Value = 3503.5
$ws.Range("K67").Value = 2000.5
$ws.Range("L67").Value = 3503.5
$ws.Range("M67").Value = -1220.5
$ws.Range("N67").Value = -5063.5
# Row 134
$ws.Range("H134").Value = 7236.4
$ws.Range("I134").Value = 4507.5386
$ws.Range("K134").Value = 13522.6158
$ws.Range("M134").Value = -10987.6158

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1596.2273
$ws.Range("I31").Value = 902.1111
$ws.Range("J31").Value = 2076.7693
$ws.Range("K31").Value = 902.1111
$ws.Range("L31").Value = 2076.7693
$ws.Range("M31").Value = -607.1111
$ws.Range("N31").Value = -2666.7693
# Row 34
$ws.Range("H34").Value = 1596.2273
$ws.Range("I34").Value = 902.1111
$ws.Range("J34").Value = 2076.7693
$ws.Range("K34").Value = 902.1111
$ws.Range("L34").Value = 2076.7693
$ws.Range("M34").Value = -700.1111
$ws.Range("N34").Value = -2480.7693
# Row 58
$ws.Range("H58").Value = 14284.182
$ws.Range("I58").Value = 5578
$ws.Range("J58").Value = 29520
$ws.Range("K58").Value = 5578
$ws.Range("L58").Value = 29520
$ws.Range("M58").Value = -5375
$ws.Range("N58").Value = -29926
# Row 132
$ws.Range("H132").Value = 21784.508
$ws.Range("I132").Value = 15295.286
$ws.Range("K132").Value = 45885.858
$ws.Range("M132").Value = -43355.858
# Row 134
$ws.Range("H134").Value = 9118.464
$ws.Range("I134").Value = 6405.0303
$ws.Range("K134").Value = 19215.0909
$ws.Range("M134").Value = -16680.0909
# Row 136
$ws.Range("H136").Value = 14284.182
$ws.Range("I136").Value = 5578
$ws.Range("J136").Value = 29520
$ws.Range("K136").Value = 16734
$ws.Range("L136").Value = 88560
$ws.Range("M136").Value = -14184
$ws.Range("N136").Value = -93660

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 380.93332
$ws.Range("I2").Value = 122.5
$ws.Range("J2").Value = 676.2857
$ws.Range("K2").Value = 735
$ws.Range("L2").Value = 4057.7142
$ws.Range("M2").Value = -622
$ws.Range("N2").Value = -4283.7142
# Row 5
$ws.Range("H5").Value = 1780.35
$ws.Range("I5").Value = 1132.8695
$ws.Range("J5").Value = 2656.353
$ws.Range("K5").Value = 3398.6085
$ws.Range("L5").Value = 7969.059
$ws.Range("M5").Value = -3286.6085
$ws.Range("N5").Value = -8193.059000000001
# Row 14
$ws.Range("H14").Value = 4516.5864
$ws.Range("I14").Value = 4516.5864
$ws.Range("K14").Value = 13549.7592
$ws.Range("M14").Value = -13376.7592
# Row 80
$ws.Range("H80").Value = 25825.572
$ws.Range("I80").Value = 100000
$ws.Range("J80").Value = 13463.167
$ws.Range("K80").Value = 300000
$ws.Range("L80").Value = 40389.501
$ws.Range("M80").Value = -299064
$ws.Range("N80").Value = -42261.501
# Row 83
$ws.Range("H83").Value = 25825.572
$ws.Range("I83").Value = 100000
$ws.Range("J83").Value = 13463.167
$ws.Range("K83").Value = 900000
$ws.Range("L83").Value = 121168.503
$ws.Range("M83").Value = -895320
$ws.Range("N83").Value = -130528.503
# Row 135
$ws.Range("H135").Value = 1780.35
$ws.Range("I135").Value = 1132.8695
$ws.Range("J135").Value = 2656.353
$ws.Range("K135").Value = 10195.8255
$ws.Range("L135").Value = 23907.177
$ws.Range("M135").Value = -7660.825500000001
$ws.Range("N135").Value = -28977.177

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 16141.5
$ws.Range("J132").Value = 49998
$ws.Range("L132").Value = 149994
$ws.Range("N132").Value = -155054

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2499.6667
$ws.Range("I7").Value = 2249.5
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 2249.5
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -2137.5
$ws.Range("N7").Value = -3224
# Row 55
$ws.Range("H55").Value = 752.4
$ws.Range("I55").Value = 983.5714
$ws.Range("K55").Value = 983.5714
$ws.Range("M55").Value = -810.5714
# Row 61
$ws.Range("H61").Value = 6245.5
$ws.Range("I61").Value = 4994.5
$ws.Range("J61").Value = 7496.5
$ws.Range("K61").Value = 4994.5
$ws.Range("L61").Value = 7496.5
$ws.Range("M61").Value = -4792.5
$ws.Range("N61").Value = -7900.5
# Row 113
$ws.Range("H113").Value = 6245.5
$ws.Range("I113").Value = 4994.5
$ws.Range("J113").Value = 7496.5
$ws.Range("K113").Value = 4994.5
$ws.Range("L113").Value = 7496.5
$ws.Range("M113").Value = -2824.5
$ws.Range("N113").Value = -11836.5
# Row 126
$ws.Range("H126").Value = 2499.6667
$ws.Range("I126").Value = 2249.5
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 6748.5
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -4278.5
$ws.Range("N126").Value = -13940
# Row 132
$ws.Range("H132").Value = 3009720
$ws.Range("I132").Value = 3589246.2
$ws.Range("J132").Value = 15500.667
$ws.Range("K132").Value = 10767738.6
$ws.Range("L132").Value = 46502.001
$ws.Range("M132").Value = -10765208.6
$ws.Range("N132").Value = -51562.001
# Row 136
$ws.Range("H136").Value = 5054694.5
$ws.Range("I136").Value = 6176543.5
$ws.Range("J136").Value = 6375
$ws.Range("K136").Value = 18529630.5
$ws.Range("L136").Value = 19125
$ws.Range("M136").Value = -18527080.5
$ws.Range("N136").Value = -24225

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 63839.832
$ws.Range("I126").Value = 8882.375
$ws.Range("K126").Value = 26647.125
$ws.Range("M126").Value = -24177.125
# Row 132
$ws.Range("H132").Value = 19535.807
$ws.Range("I132").Value = 17792
$ws.Range("J132").Value = 22706.363
$ws.Range("K132").Value = 53376
$ws.Range("L132").Value = 68119.08900000001
$ws.Range("M132").Value = -50846
$ws.Range("N132").Value = -73179.08900000001
